$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1315.6666
$ws.Cells.Item(32, 9).Value = 1473.5
$ws.Cells.Item(32, 11).Value = 1473.5
$ws.Cells.Item(32, 13).Value = -1147.5

$ws.Cells.Item(33, 8).Value = 911.375
$ws.Cells.Item(33, 9).Value = 911.375
$ws.Cells.Item(33, 11).Value = 911.375
$ws.Cells.Item(33, 13).Value = -682.375

$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(34, 8).Value = 3122
$ws.Cells.Item(34, 9).Value = 3122
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 3122
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).Value = -2919

$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(36, 8).Value = 3122
$ws.Cells.Item(36, 9).Value = 3122
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 3122
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 14).Value = -2407

$ws.Cells.Item(40, 8).Value = 3835.3572
$ws.Cells.Item(40, 9).Value = 2743.889
$ws.Cells.Item(40, 10).Value = 5800
$ws.Cells.Item(40, 11).Value = 2743.889
$ws.Cells.Item(40, 12).Value = 5800
$ws.Cells.Item(40, 13).Value = -2568.889
$ws.Cells.Item(40, 14).Value = -6150

$ws.Cells.Item(76, 8).Value = 1500
$ws.Cells.Item(76, 9).Value = 1500
$ws.Cells.Item(76, 11).Value = 1500
$ws.Cells.Item(76, 13).Value = -1185

$ws.Cells.Item(79, 8).Value = 1500
$ws.Cells.Item(79, 9).Value = 1500
$ws.Cells.Item(79, 11).Value = 1500
$ws.Cells.Item(79, 13).Value = -408

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7946.3335
$ws.Cells.Item(32, 9).Value = 7946.3335
$ws.Cells.Item(32, 11).Value = 7946.3335
$ws.Cells.Item(32, 13).Value = -7659.3335

$ws.Cells.Item(36, 8).Value = 1571.4286
$ws.Cells.Item(36, 9).Value = 1350
$ws.Cells.Item(36, 11).Value = 1350
$ws.Cells.Item(36, 13).Value = -1004

$ws.Cells.Item(61, 8).Value = 1841
$ws.Cells.Item(61, 9).Value = 1841
$ws.Cells.Item(61, 11).Value = 1841
$ws.Cells.Item(61, 13).Value = -1629

$ws.Cells.Item(74, 8).Value = 7040.3335
$ws.Cells.Item(74, 9).Value = 6771.364
$ws.Cells.Item(74, 11).Value = 6771.364
$ws.Cells.Item(74, 13).Value = -5897.364

$ws.Cells.Item(77, 8).Value = 7040.3335
$ws.Cells.Item(77, 9).Value = 6771.364
$ws.Cells.Item(77, 11).Value = 33856.82
$ws.Cells.Item(77, 13).Value = -29488.82

$ws.Cells.Item(110, 8).Value = 3026.5
$ws.Cells.Item(110, 9).Value = 1690.8462
$ws.Cells.Item(110, 10).Value = 8814.333000000001
$ws.Cells.Item(110, 11).Value = 1690.8462
$ws.Cells.Item(110, 12).Value = 8814.333000000001
$ws.Cells.Item(110, 13).Value = 354.1538
$ws.Cells.Item(110, 14).Value = -12904.333

$ws.Cells.Item(122, 8).Value = 3262.3125
$ws.Cells.Item(122, 9).Value = 3209
$ws.Cells.Item(122, 10).Value = 3379.6
$ws.Cells.Item(122, 11).Value = 9627
$ws.Cells.Item(122, 12).Value = 10138.8
$ws.Cells.Item(122, 13).Value = -7177
$ws.Cells.Item(122, 14).Value = -15038.8

$ws.Cells.Item(136, 8).Value = 1841
$ws.Cells.Item(136, 9).Value = 1841
$ws.Cells.Item(136, 11).Value = 5523
$ws.Cells.Item(136, 13).Value = -2973

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 3181.8
$ws.Cells.Item(94, 10).Value = 4401.75
$ws.Cells.Item(94, 12).Value = 4401.75
$ws.Cells.Item(94, 14).Value = -5303.75

$ws.Cells.Item(134, 8).Value = 2069.9092
$ws.Cells.Item(134, 9).Value = 1996.9
$ws.Cells.Item(134, 11).Value = 5990.700000000001
$ws.Cells.Item(134, 13).Value = -3455.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 852.75
$ws.Cells.Item(16, 9).Value = 799.3333
$ws.Cells.Item(16, 10).Value = 1013
$ws.Cells.Item(16, 11).Value = 799.3333
$ws.Cells.Item(16, 12).Value = 1013
$ws.Cells.Item(16, 13).Value = -512.3333
$ws.Cells.Item(16, 14).Value = -1587

$ws.Cells.Item(105, 8).Value = 5232.5
$ws.Cells.Item(105, 9).Value = 7250
$ws.Cells.Item(105, 10).Value = 4223.75
$ws.Cells.Item(105, 11).Value = 7250
$ws.Cells.Item(105, 12).Value = 4223.75
$ws.Cells.Item(105, 13).Value = -5503
$ws.Cells.Item(105, 14).Value = -7717.75

$ws.Cells.Item(113, 8).Value = 852.75
$ws.Cells.Item(113, 9).Value = 799.3333
$ws.Cells.Item(113, 10).Value = 1013
$ws.Cells.Item(113, 11).Value = 799.3333
$ws.Cells.Item(113, 12).Value = 1013
$ws.Cells.Item(113, 13).Value = 1370.6667
$ws.Cells.Item(113, 14).Value = -5353

$ws.Cells.Item(141, 8).Value = 99993.5
$ws.Cells.Item(141, 10).Value = 99993.5
$ws.Cells.Item(141, 12).Value = 99993.5
$ws.Cells.Item(141, 14).Value = -110353.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 364.25
$ws.Cells.Item(12, 10).Value = 285.66666
$ws.Cells.Item(12, 12).Value = 856.9999799999999
$ws.Cells.Item(12, 14).Value = -1202.99998

$ws.Cells.Item(113, 8).Value = 1257.2858
$ws.Cells.Item(113, 10).Value = 1257.2858
$ws.Cells.Item(113, 12).Value = 3771.8574
$ws.Cells.Item(113, 14).Value = -8111.857400000001

$ws.Cells.Item(122, 8).Value = 3002.5454
$ws.Cells.Item(122, 10).Value = 3303.6316
$ws.Cells.Item(122, 12).Value = 29732.6844
$ws.Cells.Item(122, 14).Value = -34632.6844

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3801.6667
$ws.Cells.Item(113, 9).Value = 3801.6667
$ws.Cells.Item(113, 11).Value = 3801.6667
$ws.Cells.Item(113, 13).Value = -1631.6667

$ws.Cells.Item(132, 8).Value = 3419.0588
$ws.Cells.Item(132, 9).Value = 3217.4614
$ws.Cells.Item(132, 11).Value = 9652.3842
$ws.Cells.Item(132, 13).Value = -7122.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 5486.0527
$ws.Cells.Item(22, 10).Value = 7687.5
$ws.Cells.Item(22, 12).Value = 7687.5
$ws.Cells.Item(22, 14).Value = -8277.5

$ws.Cells.Item(27, 8).Value = 5486.0527
$ws.Cells.Item(27, 10).Value = 7687.5
$ws.Cells.Item(27, 12).Value = 7687.5
$ws.Cells.Item(27, 14).Value = -7901.5

$ws.Cells.Item(53, 13).ClearContents()
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 11).Value = 0

$ws.Cells.Item(55, 8).Value = 2105.9285
$ws.Cells.Item(55, 9).Value = 5312.5
$ws.Cells.Item(55, 10).Value = 823.3
$ws.Cells.Item(55, 11).Value = 5312.5
$ws.Cells.Item(55, 12).Value = 823.3
$ws.Cells.Item(55, 13).Value = -5139.5
$ws.Cells.Item(55, 14).Value = -1169.3

$ws.Cells.Item(82, 8).Value = 2644.4443
$ws.Cells.Item(82, 9).Value = 2644.4443
$ws.Cells.Item(82, 11).Value = 2644.4443
$ws.Cells.Item(82, 13).Value = -2283.4443

$ws.Cells.Item(85, 8).Value = 2644.4443
$ws.Cells.Item(85, 9).Value = 2644.4443
$ws.Cells.Item(85, 11).Value = 2644.4443
$ws.Cells.Item(85, 13).Value = -1396.4443

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3559.1
$ws.Cells.Item(96, 9).Value = 4170.25
$ws.Cells.Item(96, 10).Value = 1114.5
$ws.Cells.Item(96, 11).Value = 4170.25
$ws.Cells.Item(96, 12).Value = 1114.5
$ws.Cells.Item(96, 13).Value = -2797.25
$ws.Cells.Item(96, 14).Value = -3860.5

$ws.Cells.Item(136, 8).Value = 5087.125
$ws.Cells.Item(136, 9).Value = 5087.125
$ws.Cells.Item(136, 11).Value = 15261.375
$ws.Cells.Item(136, 13).Value = -12711.375
